$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the "last updated" timestamp string in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 19:52"

# --- 2) Direct numeric updates (new COVID-19 figures) for existing rows ---
# Row 4 (Estados Unidos)
$ws.Range("B4").Value = 903298
$ws.Range("C4").Value = 16856
$ws.Range("E4").Value = 762083
$ws.Range("G4").Value = 718
$ws.Range("H4").Value = 50954

# Row 7 (Francia)
$ws.Range("B7").Value = 159828
$ws.Range("C7").Value = 1645
$ws.Range("D7").Value = 43493
$ws.Range("E7").Value = 94090
$ws.Range("F7").Value = 4870
$ws.Range("G7").Value = 389
$ws.Range("H7").Value = 22245

# Row 8 (Alemania)
$ws.Range("B8").Value = 154159
$ws.Range("C8").Value = 1030
$ws.Range("E8").Value = 41706
$ws.Range("G8").Value = 78
$ws.Range("H8").Value = 5653

# Row 16 (Iran)
$ws.Range("B16").Value = 43551
$ws.Range("C16").Value = 1441
$ws.Range("D16").Value = 15444
$ws.Range("E16").Value = 25813
$ws.Range("G16").Value = 147
$ws.Range("H16").Value = 2294

# Row 19 (Belgica)
$ws.Range("B19").Value = 24434
$ws.Range("C19").Value = 1395
$ws.Range("D19").Value = 5457
$ws.Range("E19").Value = 18197
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = 780

# Row 69
$ws.Range("D69").Value = 621
$ws.Range("E69").Value = 1149

# --- 3) Relocate "Guinea Ecuatorial" within the country list ---
# It used to sit right after "Cabo Verde" (old row 150); it now belongs
# right after "Paraguay" (row 128), pushing "Islas Feroe" and everything
# through "Cabo Verde" down by one row.
$ws.Rows("128:128").Insert()

$ws.Range("A128").Value = "Guinea Ecuatorial"
$ws.Range("B128").Value = 214
$ws.Range("C128").Value = 5
$ws.Range("D128").Value = 7
$ws.Range("E128").Value = 206
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 1

# The old "Guinea Ecuatorial" row has now shifted down to row 151 (it used
# to be row 150, and the insert above pushed it down by one); remove it so
# the table doesn't contain the country twice.
$ws.Rows("151:151").Delete()
